$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# Rename header cells: "_old" -> "_FV2310", "_new" -> "_FV2404"
for ($i = 1; $i -le 21; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $val = $cell.Value2
    if ($val -like "*_old") {
        $cell.Value2 = ($val -replace "_old$", "_FV2310")
    } elseif ($val -like "*_new") {
        $cell.Value2 = ($val -replace "_new$", "_FV2404")
    }
}

# Turn the used range into an Excel Table (ListObject)
$range = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (select the cell below the split point, then freeze)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
